$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before K (shifts K..AB to L..AC), mirroring the
# "Insert Column" command run against column K in the source sheet.
$leftWidth = $ws.Columns("J").ColumnWidth
$ws.Columns("K").Insert()

# Carry the neighbouring column's width onto the freshly inserted one.
$ws.Columns("K").ColumnWidth = $leftWidth

# New header cell: "DynamiteOccurrenceLinkLocation" (bold, text format).
$ws.Range("K1").Value = "DynamiteOccurrenceLinkLocation"
$ws.Range("K1").NumberFormat = "@"
$ws.Range("K1").Font.Bold = $true

# New data cell: "Main Menu" (text format).
$ws.Range("K2").Value = "Main Menu"
$ws.Range("K2").NumberFormat = "@"

# Leave the new column selected, matching the edited workbook's view state.
$ws.Range("K1:K2").Select() | Out-Null
